$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "54.265.53"
$ws.Range("E2").Value = "  -0.16%  "
$ws.Range("D3").Value = "2.281.40"
$ws.Range("E3").Value = "  +0.12%  "
$ws.Range("E4").Value = "  -0.01%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "498.77"
$ws.Range("E5").Value = "  +1.01%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "128.66"
$ws.Range("E6").Value = "  +1.05%  "
$ws.Range("E7").Value = "  -0.19%  "
$ws.Range("E8").Value = "  -0.07%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.0954"
$ws.Range("E9").Value = "  +2.29%  "
$ws.Range("E10").Value = "  +1.15%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.333"
$ws.Range("E11").Value = "  +2.89%  "
$ws.Range("E12").Value = "  +1.28%  "
$ws.Range("D13").Value = "2.684.78"
$ws.Range("E13").Value = "  +0.23%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "22.70"
$ws.Range("E14").Value = "  +5.50%  "
$ws.Range("D15").Value = "54.217.61"
$ws.Range("E15").Value = "  -0.10%  "
$ws.Range("E16").Value = "  +0.01%  "
$ws.Range("D17").Value = "2.301.04"
$ws.Range("E17").Value = "  +1.55%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "10.24"
$ws.Range("E18").Value = "  +4.18%  "
$ws.Range("E19").Value = "  +1.76%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "304.54"
$ws.Range("E20").Value = "  +1.92%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "6.41"
$ws.Range("E21").Value = "  +2.16%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.999"
$ws.Range("E22").Value = "  -0.14%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "61.84"
$ws.Range("E23").Value = "  -2.84%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "1.00"
$ws.Range("E24").Value = "  -0.15%  "
$ws.Range("E25").Value = "  +2.37%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "7.31"
$ws.Range("E26").Value = "  +2.35%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "174.77"
$ws.Range("E27").Value = "  +7.16%  "
$ws.Range("E28").Value = "  +0.64%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "5.96"
$ws.Range("E29").Value = "  +2.50%  "
$ws.Range("D30").Value = "0.0₃0686"
$ws.Range("E30").Value = "  +0.60%  "
$ws.Range("E31").Value = "  +0.85%  "
$ws.Range("E32").Value = "  -0.03%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "17.78"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.996"
$ws.Range("E34").Value = "  -0.28%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.934"
$ws.Range("E35").Value = "  +9.42%  "
$ws.Range("E36").Value = "  +0.85%  "
$ws.Range("E37").Value = "  +2.59%  "
$ws.Range("B38").Value = "OKB"
$ws.Range("C38").Value = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "36.06"
$ws.Range("E38").Value = "  +1.52%  "
$ws.Range("B39").Value = "PolygonEcosystemToken"
$ws.Range("C39").Value = "https://coinranking.com/coin/iDZ0tG-wI+polygonecosystemtoken-pol"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.374"
$ws.Range("E39").Value = "  -0.33%  "
$ws.Range("B40").Value = "Stacks"
$ws.Range("C40").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "1.41"
$ws.Range("E40").Value = "  +0.50%  "
$ws.Range("B41").Value = "Filecoin"
$ws.Range("C41").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "3.38"
$ws.Range("E41").Value = "  +1.48%  "
$ws.Range("B42").Value = "Aave"
$ws.Range("C42").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "124.93"
$ws.Range("E42").Value = "  -1.26%  "
$ws.Range("B43").Value = "RenderToken"
$ws.Range("C43").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "4.76"
$ws.Range("E43").Value = "  -1.47%  "
$ws.Range("B44").Value = "Hedera"
$ws.Range("C44").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.0493"
$ws.Range("E44").Value = "  +2.90%  "
$ws.Range("B45").Value = "Stellar"
$ws.Range("C45").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.0895"
$ws.Range("E45").Value = "  +0.52%  "
$ws.Range("B46").Value = "Mantle"
$ws.Range("C46").Value = "https://coinranking.com/coin/BoI4ux0nd+mantle-mnt"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.546"
$ws.Range("E46").Value = "  +0.06%  "
$ws.Range("B47").Value = "Bittensor"
$ws.Range("C47").Value = "https://coinranking.com/coin/pgv7xSFi6+bittensor-tao"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "240.06"
$ws.Range("E47").Value = "  -0.63%  "
$ws.Range("B48").Value = "Polygon"
$ws.Range("C48").Value = "https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.373"
$ws.Range("E48").Value = "  -0.33%  "
$ws.Range("B49").Value = "VeChain"
$ws.Range("C49").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.0206"
$ws.Range("E49").Value = "  +1.32%  "
$ws.Range("B50").Value = "WhiteBITCoin"
$ws.Range("C50").Value = "https://coinranking.com/coin/GE4c3_TbB+whitebitcoin-wbt"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "10.77"
$ws.Range("E50").Value = "  +1.00%  "
$ws.Range("B51").Value = "InjectiveProtocol"
$ws.Range("C51").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "16.27"
$ws.Range("E51").Value = "  +0.12%  "
